$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above row 672 (shifts old rows 672:682 down to 679:689)
$ws.Range("A672:A678").EntireRow.Insert()

# Common column values shared by every row in this data block
$colA = 7
$colB = "Terminal Hortofrutícola Agro Chillán"
$colC = "Ñuble"
$colE = 16
$colF = "Fruta"
$colG = 100104
$colH = "Frutos de pepita"
$colI = 100104002
$colJ = "Manzana"
$colQ = "`$/caja 16 kilos empedrada"
$colR = "Provincia de Curicó"
$colT = 16

function Set-PriceRow($r, $d, $k, $l, $m, $n, $o, $p, $s) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $colR
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $colT
}

Set-PriceRow 672 44656 "Fuji royal"    "Especial" 60  10000 10000 10000 625
Set-PriceRow 673 44656 "Fuji royal"    "Primera"  120 8500  9000  8750  547
Set-PriceRow 674 44656 "Fuji royal"    "Segunda"  60  8000  8000  8000  500
Set-PriceRow 675 44656 "Granny Smith"  "Especial" 60  10000 10000 10000 625
Set-PriceRow 676 44656 "Granny Smith"  "Primera"  120 8500  9000  8750  547
Set-PriceRow 677 44656 "Granny Smith"  "Segunda"  60  8000  8000  8000  500
Set-PriceRow 678 44656 "Royal Gala"    "Primera"  120 7500  8000  7750  484
